# Fix the ordering of names/emails in the "Recorded By" column (G).
# Two specific value patterns need their comma-separated parts reordered:
#   "backup@backdoor.com, System, system"  ->  "system, backup@backdoor.com, System"
#   "System, dnasr281@gmail.com"           ->  "dnasr281@gmail.com, System"
# This mirrors swaps that look like someone moving "system"/"System" entries
# to reflect a corrected recorded-by ordering across the whole sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$rowCount = $used.Rows.Count

for ($r = 1; $r -le $rowCount; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G
    $val = $cell.Value2

    if ($val -eq "backup@backdoor.com, System, system") {
        $cell.Value = "system, backup@backdoor.com, System"
    }
    elseif ($val -eq "System, dnasr281@gmail.com") {
        $cell.Value = "dnasr281@gmail.com, System"
    }
}
